# Update as of 2024-02-24
# Adds new Ingreso (contribution) rows, a new Gastos (expense) row, and a
# new Cuentas por cobrar (receivable) row to the aportes workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Ingreso": 14 new contribution rows for date 2024-02-24 (45347)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Ingreso")

# Carry down the existing row formatting so no new cell styles are minted:
#   rows 649-660 look like row 647 (plain number style on column C)
#   rows 661-662 look like row 648 (right-aligned / s=2 style on column C)
$ws1.Range("A647:D647").Copy()
$ws1.Range("A649:D660").PasteSpecial(-4122)
$ws1.Range("A648:D648").Copy()
$ws1.Range("A661:D662").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$fecha = 45347
$ingresoRows = @(
    @("Julio",    200),
    @("Omaury",   100),
    @("Elio",     100),
    @("Yeyo",     145),
    @("Randy",    100),
    @("Joel",     100),
    @("Invitados",100),
    @("Alfredo",  100),
    @("Kibelo",     0),
    @("Rayder",     0),
    @("Johan",      0),
    @("Frandy",     0),
    @("Javier",     0),
    @("Rubio",      0)
)

$row = 649
foreach ($entry in $ingresoRows) {
    $ws1.Cells.Item($row, 1).Value = $fecha
    $ws1.Cells.Item($row, 2).Value = $entry[0]
    $ws1.Cells.Item($row, 3).Value = $entry[1]
    $ws1.Cells.Item($row, 4).Value = "Aporte"
    $row = $row + 1
}

# ---------------------------------------------------------------------
# Sheet "Gastos": new expense row for 2024-02-24
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Gastos")

$ws2.Range("A81:C81").Copy()
$ws2.Range("A82:C82").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws2.Cells.Item(82, 1).Value = 45347
$ws2.Cells.Item(82, 2).Value = "Arbitro y agua"
$ws2.Cells.Item(82, 3).Formula = "=150+800"

# ---------------------------------------------------------------------
# Sheet "Cuentas por cobrar": new receivable row
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Cuentas por cobrar")

$ws3.Range("A7:F7").Copy()
$ws3.Range("A8:F8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws3.Cells.Item(8, 1).Value = 45347
$ws3.Cells.Item(8, 2).Value = "Jordan"
$ws3.Cells.Item(8, 3).Value = "Tecnica"
$ws3.Cells.Item(8, 4).Value = 100
$ws3.Cells.Item(8, 5).ClearContents()
$ws3.Cells.Item(8, 6).Value = 'Le dijo "Mamagüevo" a los que no entraban a rebotar'

# ---------------------------------------------------------------------
# Restore the view state recorded in the saved file. "Ingreso" is the
# tab left active/selected, so activate it last.
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("B82:C82").Select()

$ws3.Activate()
$ws3.Range("F5").Select()

$ws1.Activate()
$ws1.Range("C655").Select()
